$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing categories (rows 184-218) so the sheet's
# used range shrinks from A1:E218 down to A1:E183.
$ws.Range("A184:E218").EntireRow.Delete()

# Rebuild rows 2-183 (category, date, forecast) to match the refreshed
# forecast run.
$categories = @('groceries','groceries','groceries','groceries','groceries','groceries','groceries','ecom','ecom','ecom','ecom','ecom','ecom','ecom','delivery food','delivery food','delivery food','delivery food','delivery food','delivery food','delivery food','medicine/health','medicine/health','medicine/health','medicine/health','medicine/health','medicine/health','medicine/health','sport','sport','sport','sport','sport','sport','sport','subscriptions','subscriptions','subscriptions','subscriptions','subscriptions','subscriptions','subscriptions','taxi car','taxi car','taxi car','taxi car','taxi car','taxi car','taxi car','cafes','cafes','cafes','cafes','cafes','cafes','cafes','fun','fun','fun','fun','fun','fun','fun','school','school','school','school','school','school','school','massage','massage','massage','massage','massage','massage','massage','hotel accomondation and fun','hotel accomondation and fun','hotel accomondation and fun','hotel accomondation and fun','hotel accomondation and fun','hotel accomondation and fun','hotel accomondation and fun','steam&PS&Google play','steam&PS&Google play','steam&PS&Google play','steam&PS&Google play','steam&PS&Google play','steam&PS&Google play','steam&PS&Google play','education','education','education','education','education','education','education','clothes','clothes','clothes','clothes','clothes','clothes','clothes','phone, internet','phone, internet','phone, internet','phone, internet','phone, internet','phone, internet','phone, internet','rent + communal','rent + communal','rent + communal','rent + communal','rent + communal','rent + communal','rent + communal','beauty','beauty','beauty','beauty','beauty','beauty','beauty','shop (not groceries)','shop (not groceries)','shop (not groceries)','shop (not groceries)','shop (not groceries)','shop (not groceries)','shop (not groceries)','public transport','public transport','public transport','public transport','public transport','public transport','public transport','present','present','present','present','present','present','present','Ira transfer','Ira transfer','Ira transfer','Ira transfer','Ira transfer','Ira transfer','Ira transfer','car rent','car rent','car rent','car rent','car rent','car rent','car rent','parking','parking','parking','parking','parking','parking','parking','gas','gas','gas','gas','gas','gas','gas','Mira activities','Mira activities','Mira activities','Mira activities','Mira activities','Mira activities','Mira activities')
$dates = @(45767,45768,45769,45770,45771,45772,45773,45766,45767,45768,45769,45770,45771,45772,45766,45767,45768,45769,45770,45771,45772,45767,45768,45769,45770,45771,45772,45773,45765,45766,45767,45768,45769,45770,45771,45744,45745,45746,45747,45748,45749,45750,45708,45709,45710,45711,45712,45713,45714,45766,45767,45768,45769,45770,45771,45772,45767,45768,45769,45770,45771,45772,45773,45763,45764,45765,45766,45767,45768,45769,45766,45767,45768,45769,45770,45771,45772,45682,45683,45684,45685,45686,45687,45688,45684,45685,45686,45687,45688,45689,45690,45755,45756,45757,45758,45759,45760,45761,45766,45767,45768,45769,45770,45771,45772,45766,45767,45768,45769,45770,45771,45772,45766,45767,45768,45769,45770,45771,45772,45761,45762,45763,45764,45765,45766,45767,45764,45765,45766,45767,45768,45769,45770,45740,45741,45742,45743,45744,45745,45746,45756,45757,45758,45759,45760,45761,45762,45766,45767,45768,45769,45770,45771,45772,45755,45756,45757,45758,45759,45760,45761,45755,45756,45757,45758,45759,45760,45761,45767,45768,45769,45770,45771,45772,45773,45767,45768,45769,45770,45771,45772,45773)
$forecasts = @(26.075,26.075,26.075,26.075,26.075,26.075,26.075,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,3.586,3.586,3.586,3.586,3.586,3.586,3.586,11.86833333333333,11.86833333333333,11.86833333333333,11.86833333333333,11.86833333333333,11.86833333333333,11.86833333333333,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1.666666666666667,1.666666666666667,1.666666666666667,1.666666666666667,1.666666666666667,1.666666666666667,1.666666666666667,37,37,37,37,37,37,37,18.8,18.8,18.8,18.8,18.8,18.8,18.8,100,100,100,100,100,100,100,0,0,0,0,0,0,0,0.4,0.4,0.4,0.4,0.4,0.4,0.4,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,75,75,75,75,75,75,75,0,0,0,0,0,0,0,6.666666666666667,6.666666666666667,6.666666666666667,6.666666666666667,6.666666666666667,6.666666666666667,6.666666666666667,0,0,0,0,0,0,0,4,4,4,4,4,4,4,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $categories[$i]
    $ws.Cells.Item($r, 2).Value = $dates[$i]
    $ws.Cells.Item($r, 3).Value = $forecasts[$i]
}
